$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B32").Value = 66452
$ws.Range("F32").Value = 69
$ws.Range("G32").Value = 2121.06
$ws.Range("B33").Value = 51755
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 30.74
$ws.Range("F42").Value = 152
$ws.Range("G42").Value = 29909.04
$ws.Range("F44").Value = 574
$ws.Range("G44").Value = 20899.34
$ws.Range("F47").Value = 264
$ws.Range("G47").Value = 50922.96
$ws.Range("F54").Value = 269
$ws.Range("G54").Value = 15090.9
$ws.Range("F59").Value = 146
$ws.Range("G59").Value = 8614
$ws.Range("F65").Value = 129
$ws.Range("G65").Value = 10052.97
$ws.Range("F68").Value = 282
$ws.Range("G68").Value = 73525.86
$ws.Range("B73").Value = 276505.26
$ws.Range("F94").Value = 1
$ws.Range("G94").Value = 2169.54
$ws.Range("B101").Value = 112435.89
$ws.Range("F148").Value = 2
$ws.Range("G148").Value = 55979.22
$ws.Range("B152").Value = 206389.92
$ws.Range("F187").Value = 65
$ws.Range("G187").Value = 8668.4
$ws.Range("B191").Value = 48126.68
$ws.Range("F233").Value = 65
$ws.Range("G233").Value = 5820.1
$ws.Range("B249").Value = 98998.25
$ws.Range("F275").Value = 30
$ws.Range("G275").Value = 3060.3
$ws.Range("F278").Value = 141
$ws.Range("G278").Value = 14869.86
$ws.Range("F279").Value = 228
$ws.Range("G279").Value = 25602.12
$ws.Range("B282").Value = 111318.24
$ws.Range("F286").Value = 2169
$ws.Range("G286").Value = 40126.5
$ws.Range("B293").Value = 61804.5
$ws.Range("F315").Value = 11
$ws.Range("G315").Value = 2526.48
$ws.Range("B317").Value = 22187.95
$ws.Range("B324").Value = 66188
$ws.Range("C324").Value = 'HIM-Baby Care Gift Pack (Ww)1'
$ws.Range("D324").Value = 315.8
$ws.Range("E324").Value = 377.31
$ws.Range("F324").Value = 45
$ws.Range("G324").Value = 14211
$ws.Range("B325").Value = 48719
$ws.Range("C325").Value = 'HIM-BABY CARE GIFT PACK (WW)1'
$ws.Range("D325").Value = 295.75
$ws.Range("E325").Value = 353.35
$ws.Range("F325").Value = -82
$ws.Range("G325").Value = -24251.5
$ws.Range("B369").Value = 64983
$ws.Range("C369").Value = 'HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S'
$ws.Range("F369").Value = 6
$ws.Range("G369").Value = 514.08
$ws.Range("B370").Value = 66194
$ws.Range("C370").Value = 'HIM-Total Care Baby Pants Diapers-M-9s'
$ws.Range("F370").Value = 37
$ws.Range("G370").Value = 3170.16
$ws.Range("F455").Value = 65
$ws.Range("G455").Value = 14964.3
$ws.Range("F456").Value = 66
$ws.Range("G456").Value = 15846.6
$ws.Range("B460").Value = 120223.24
$ws.Range("F470").Value = 324
$ws.Range("G470").Value = 53780.76
$ws.Range("B472").Value = 95135.61
$ws.Range("B482").Value = 58047
$ws.Range("D482").Value = 105.54
$ws.Range("E482").Value = 126.1
$ws.Range("F482").Value = 34
$ws.Range("G482").Value = 3588.36
$ws.Range("B483").Value = 47097
$ws.Range("D483").Value = 112.28
$ws.Range("E483").Value = 134.16
$ws.Range("F483").Value = 15
$ws.Range("G483").Value = 1684.2
$ws.Range("F502").Value = 138
$ws.Range("G502").Value = 25000.08
$ws.Range("F516").Value = 136
$ws.Range("G516").Value = 7829.52
$ws.Range("B521").Value = 204285.5
$ws.Range("B571").Value = 53595
$ws.Range("E571").Value = 17.61
$ws.Range("F571").Value = -335
$ws.Range("G571").Value = -4934.55
$ws.Range("B572").Value = 65067
$ws.Range("E572").Value = 15.65
$ws.Range("F572").Value = 126
$ws.Range("G572").Value = 1855.98
$ws.Range("F579").Value = 72
$ws.Range("G579").Value = 1936.08
$ws.Range("B586").Value = 40334.13
$ws.Range("F606").Value = 38
$ws.Range("G606").Value = 5368.26
$ws.Range("B614").Value = 143359.89
$ws.Range("B642").Value = 53319
$ws.Range("E642").Value = 310.64
$ws.Range("F642").Value = -6
$ws.Range("G642").Value = -1643.52
$ws.Range("B643").Value = 64810
$ws.Range("E643").Value = 291.22
$ws.Range("F643").Value = 2
$ws.Range("G643").Value = 547.84
$ws.Range("B661").Value = 60025
$ws.Range("E661").Value = 37.22
$ws.Range("F661").Value = -98
$ws.Range("G661").Value = -3217.34
$ws.Range("B662").Value = 64833
$ws.Range("E662").Value = 34.9
$ws.Range("F662").Value = 89
$ws.Range("G662").Value = 2921.87
$ws.Range("F687").Value = 544
$ws.Range("G687").Value = 29854.72
$ws.Range("B694").Value = 182042.74
$ws.Range("F702").Value = 63
$ws.Range("G702").Value = 2720.34
$ws.Range("F704").Value = 59
$ws.Range("G704").Value = 2547.62
$ws.Range("B707").Value = 40056.84
$ws.Range("F719").Value = 83
$ws.Range("G719").Value = 5137.7
$ws.Range("B731").Value = 37989.52
$ws.Range("F734").Value = 322
$ws.Range("G734").Value = 39235.7
$ws.Range("B741").Value = 48961.15
$ws.Range("F798").Value = 259
$ws.Range("G798").Value = 34472.9
$ws.Range("B801").Value = 36213.73
$ws.Range("F806").Value = 138
$ws.Range("G806").Value = 15015.78
$ws.Range("F809").Value = 30
$ws.Range("G809").Value = 2556
$ws.Range("F810").Value = 289
$ws.Range("G810").Value = 43468.49
$ws.Range("F811").Value = 32
$ws.Range("G811").Value = 4682.56
$ws.Range("F813").Value = 24
$ws.Range("G813").Value = 667.92
$ws.Range("F815").Value = 151
$ws.Range("G815").Value = 23355.17
$ws.Range("F816").Value = 74
$ws.Range("G816").Value = 10564.98
$ws.Range("F836").Value = 15
$ws.Range("G836").Value = 854.85
$ws.Range("B837").Value = 318707.05
$ws.Range("F877").Value = 97
$ws.Range("G877").Value = 7790.07
$ws.Range("B884").Value = 24045.78
$ws.Range("F917").Value = 17
$ws.Range("G917").Value = 7423.9
$ws.Range("F919").Value = 29
$ws.Range("G919").Value = 4916.37
$ws.Range("F932").Value = 17
$ws.Range("G932").Value = 3639.02
$ws.Range("B935").Value = 110315.8
$ws.Range("B941").Value = 4869558.01
$ws.Range("B942").Value = 4869558.01
